$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet "My Series" -> "Data"
$ws.Name = "Data"

# A11 label: "Function Description" -> "Function Information"
$ws.Range("A11").Value = "Function Information"

# B21 Kurtosis value tweak
$ws.Range("B21").Value = 0.2499825759175085

# Update the CDM add-in comment text stored on A1
$comment = $ws.Range("A1").Comment
$comment.Text('Jx0AAB+LCAAAAAAAAAOlWVtvI0kV/istP4GE0217MptElV75lqyFHUe2Qyb7gsrdlbhIX0xXdRK/LRJo0bIIITSLlqt4WoTEMIJdaZnh8l9Wk8zwxF/g1KVvtrPjDqNo0nXOd6pOnTq3qqB3b3zPuCIRo2GwX6ltWRWDBE7o0uBivxLz82rtceVdG3VvHOId4wj7hAPYAKmA7d0wul+ZcT7fM83r6+ut68ZWGF2YdcuqmU8G/bEzIz6u0oBxHDikkkq5b5eq2Kjt+gPCsYs5VpL7ld64t9Um1OkAbYADfEGirVbMaEAY6wacckqYkIwI5qTdGXxHbcyubz3eqiFzhZ4hWzH1XIUrIBVd42BZMqE+setWbadq7VatxsR6vGfV92ogt9N4PxFMgaiPGR+T6Io6kjDm2J9LcWvXaliPrXrtMTLXgmCuzAA2GnruiFxRRtw28TxWyiKmPsCmw2HX5YxpITMnqyd6uAqHEZ7PJpR7pJwao0HL8AOtSzaJjQ7CiDhgvwepdESuh5E262TeB+5kRiO+6OBF6blOGImGc2GkcqI26oQBb3ok4idzOGvigisAw+ZRTJB5DzMT6lDmwDcNYuLa59hjeaECE52G0SWbY4ccQRybYo7rwAuxCw7HKePUyRZdYaDjKJzDjLB4K/TcA5hVg9cw0pl7AZhYLNsKw8tMu3VMJE9Vni+cqY95Al+ho/EsvB4G3mIcT5kT0SlxO60EvZaHREBq6XbMeOiDFhkJKVqOsoB/EIDLZNQhDvWxd+yBEZndgFkKBNSMeXhOeTv0Yj9giU5LVHQKO5qQm3SH6RgN4XADYfQw6AUJXpl5LasoMAqv0zVXGdIIOXKTOclxrzKWwR2gJce3ypEnInZ5QD2oD/mzyFGLXjGeEcLXuoTiIJEKD0TFsVuLo9ifQnhNIcau5KoMmRkfgZ+Cr4NetgVVpCp/Jpa1J39Aj5SNuoF7Py5hIlgut5ZdA94SCcGevJaHg0ugnlI+O2ome1nDQcoC9+JXeQgid+7hhSSnVsrTUC9wvNglKiH0gnPpokI3daj3stEKqQ8xbiMcLCaLOeRlRvc4fOxXoFLvMR5BL1CxnTAOeLQQmQOZGvo2GRZPA7kA9jaWOY/I92NoQRYHceC0Q3fz1VxlnZOA8s01DONIpcPNRaT1RGaMWYeIHCOT/sbyTpk9sagU3A+IHwbU2dzaYGShvfuAjbAkqjaWICq+NsZ7UNdV2ROxvrFYBP0jFLpSyzQZCx0qnVWHh5uTN+8JmQ45x7EHvRuHEnuR5t5lMmqyy2VMnoROIi/JgLbojBm0xo7rbznQPIj2b8sJfUEwoSM9HSMzjxcdkEO6wUUfBxcx9BhpXlmmp/lX1MdJhAMmtpO2FEupeD0IJXlKtTq2Sl7DWDqCSl4hcJG5hEMT4s/DCHsDMAw90G6n+yVoRQaYz/QIaptHnMTIZiaaShU1SxR/G0wWKbUNEfA6TS4RJUjsRTXhGSajIbHLAYSl18YenUYqqyalfB0PDixrDpP8KzZXslFMzgAuYlB9v00WojvPBpouXbaWMJQDi0Rqj0ePdurbVqMOnY0YI7njEcGe0YVg5sToBVeEcR/E9owRYdSFL4q9PeM9MiUUiqA0kS5DpaXzcuggyfNSlSboW6QUAdBvXFAoI6vAlJMJ2GcER94iB1Rb7YcO4O5+/O/b3zx/9eKzu4+evvnih//9x69e/fNnt89+BB93f/3b7ce/VNtUYDTBU49IhSatnR2r8Qj8LCUhYVxTtsZu7HBJOzuTHXE6RvpiJwftbq992G/JfJISE3FVUkxxZ1yEcTYcq03IheSRmoknKIg9SfKTHhe4uRJlizvcFSmi8/z7BJUtXr/87PXLP98rrQ2W9Vq13d3taq3+1lYMbsO1FVzaivULNUCAH1Wt7Wq9ngMvYdBIFYDUTj3XbtTgwl1vWLU0l7upI68DLbP0TBN8YS7JKVJbtUepC+THCVM6/gRCJGWrUMgNtIt+8ZM3f3laQGnrakpxFlBOtjFiMTMZyKmPRhNjPDwZtbvGpDsWfpLxcjg1+deA9eppPBWcKghi7H3LgOIOxcyowI2oYoTnBsHOzFhAJObisOBs66hqoQdOuazlYRTGc3UiOYGMugaZZpO1EmtyjeRJe64knYy1Bq50vf375+sE9EY6WUObvn/kaajAUaQcX0ftp/969eWHr168uHv+89svf1CYQa+TPguAn0M05Yep20PK0/VmiYJOx9KYl9Z3c/VFE8VF6jikAWe2eOfKRghEa2I2+Rv1fCh5cmJpL6AvUdB7mHVvuA5s+wiZRQLoOcdQbcPs7pkSVA7P7Pqf3/7u7tef333y/M2Hf7r96I+3H3/y+uXv3zz7g4q6u6fP7376TGf55UIgdRE3WtUEGvJ9xDFENBqidhtfffALIwi5AS2HEcuM9NUHn+YmE4rK5iSbGVq6VJGiCivQvLCQM3KqpDoU5FIR1QC0RQlrpAhdxMI5dbJF3q+KqUTcScY3epNqzIgRQjf1TdhJEZwJbyqnRVRJPX7Hqtfqmqu0EVuYYpYz/aEXTqHJSBjyAWIJUpD6eoEMK9c77A9bzX4GUUoMI5dEwg3VB0paSlFSeiwZJa6WowAXGj8n9sSb0QpslZXOnEtjpn58OW+6Iv2tf6ooIFA7jiLVEAX6LX8cz6EZTp7o7ufLZ8tc/3uketV8R5yNe50iH8Y5LhTCIlsQJF+mJs1SaarHxDuPamePhGmyIfAKT51gDv1crzqtK+grI1PknW4UhdHa5JNxEtgAOmnIKGZm8RQjz1R13W52VgkhSXjph7r56R2GHeIRXu4t28ykB+HVg2Xh7MuK9tjQc7Uxy109UrNkE+Qf9IWj/L/v+crZmlEEjZV4ACz9AJ9cXEdw3y2pjdqKFBQ3QFhdv5Af0IjxJyIT6C9FOUspZ6pDfSIuXOpDjs/sxrYiAMDMz24W1ExCl6s/q4Ren/q05LXQSuK7OAnYcj5XLVyvnKeI0nJEbqDBzM0ASXH6PSgb6h2lzGzKYSGXpvLi7ZLRixkvq9g7U0xcMrWqzpTUq49ca6e6S0ijWqvB/9ip1y1rW7x86skhc1ByXXIRMzmw7M+d9v8AhZ9z2ycdAAA=')
